# Update template.xlsx:
#  - A3 changes from an empty quoted-text cell to "DH"
#  - B3 changes from an empty quoted-text cell to "DHTC - Đơn hàng thành công"
#  - Row heights (1-3) change from 19.5 to 18.75
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "DH"
$ws.Range("B3").Value = "DHTC - Đơn hàng thành công"

$ws.Rows.Item(1).RowHeight = 18.75
$ws.Rows.Item(2).RowHeight = 18.75
$ws.Rows.Item(3).RowHeight = 18.75
